$wb = $excel.ActiveWorkbook

$wsY = $wb.Worksheets.Item("Variables Y")
$wsX = $wb.Worksheets.Item("Variables X")

# --- Sheet "Variables Y": fix indicator code c08_11 -> c18_11 ---
$wsY.Range("B3").Value = "c18_11"

# --- Sheet "Variables X": fill in preparation/readiness status column C ---
$wsX.Range("C2").Value = "listo"
$wsX.Range("C3").Value = "listo"
$wsX.Range("C4").Value = "listo"
$wsX.Range("C5").Value = "listo"
$wsX.Range("C6").Value = "no disponible"
$wsX.Range("C7").Value = "no disponible"
$wsX.Range("C8").Value = "listo"
$wsX.Range("C9").Value = "no disponible"

# Row 9 grew taller (model row with longer description)
$wsX.Rows.Item(9).RowHeight = 30

# --- restore cursor / view positions ---
$wsX.Activate()
$wsX.Application.ActiveWindow.ScrollRow = 4
$wsX.Range("C10").Select()

$wsY.Activate()
$wsY.Range("E10").Select()
